$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F32").Value = 77
$ws.Range("G32").Value = 5393.08
$ws.Range("F36").Value = 38
$ws.Range("G36").Value = 7329.82
$ws.Range("F44").Value = 24
$ws.Range("G44").Value = 437.76
$ws.Range("F49").Value = 33
$ws.Range("G49").Value = 736.5599999999999
$ws.Range("B56").Value = 45871.46
$ws.Range("F74").Value = 58
$ws.Range("G74").Value = 2308.4
$ws.Range("B75").Value = 24304.48
$ws.Range("F89").Value = 68
$ws.Range("G89").Value = 4781.76
$ws.Range("F91").Value = 117
$ws.Range("G91").Value = 7421.31
$ws.Range("F92").Value = 4
$ws.Range("G92").Value = 539.8
$ws.Range("F93").Value = 301
$ws.Range("G93").Value = 19173.7
$ws.Range("F95").Value = 97
$ws.Range("G95").Value = 7661.06
$ws.Range("F96").Value = 186
$ws.Range("G96").Value = 26105.1
$ws.Range("F104").Value = 151
$ws.Range("G104").Value = 15471.46
$ws.Range("F105").Value = 56
$ws.Range("G105").Value = 11295.76
$ws.Range("F106").Value = 32
$ws.Range("G106").Value = 4312.64
$ws.Range("F114").Value = 270
$ws.Range("G114").Value = 5254.2
$ws.Range("B115").Value = 259067.11
$ws.Range("B156").Value = 57756
$ws.Range("B157").Value = 53925
$ws.Range("F174").Value = 41
$ws.Range("G174").Value = 3362
$ws.Range("B185").Value = 19729.06
$ws.Range("F199").Value = 1
$ws.Range("G199").Value = 111.88
$ws.Range("B204").Value = 5334.6
$ws.Range("F266").Value = 46
$ws.Range("G266").Value = 1731.9
$ws.Range("F273").Value = 7
$ws.Range("G273").Value = 241.22
$ws.Range("B279").Value = 120099.49
$ws.Range("F283").Value = 210
$ws.Range("G283").Value = 21569.1
$ws.Range("F312").Value = 170
$ws.Range("G312").Value = 19033.2
$ws.Range("F323").Value = 146
$ws.Range("G323").Value = 14760.6
$ws.Range("F327").Value = 0
$ws.Range("G327").Value = 0
$ws.Range("F330").Value = 115
$ws.Range("G330").Value = 6799.95
$ws.Range("F338").Value = 13
$ws.Range("G338").Value = 6837.35
$ws.Range("F339").Value = 57
$ws.Range("G339").Value = 11386.89
$ws.Range("F340").Value = 23
$ws.Range("G340").Value = 4277.77
$ws.Range("F341").Value = 220
$ws.Range("G341").Value = 15468.2
$ws.Range("B349").Value = 375932.11
$ws.Range("B396").Value = 58047
$ws.Range("D396").Value = 105.54
$ws.Range("E396").Value = 126.1
$ws.Range("F396").Value = 62
$ws.Range("G396").Value = 6543.48
$ws.Range("B397").Value = 47097
$ws.Range("D397").Value = 112.28
$ws.Range("E397").Value = 134.16
$ws.Range("F397").Value = 15
$ws.Range("G397").Value = 1684.2
$ws.Range("F404").Value = 66
$ws.Range("G404").Value = 1622.94
$ws.Range("B407").Value = 50912.7
$ws.Range("F412").Value = 65
$ws.Range("G412").Value = 6386.25
$ws.Range("B424").Value = 48693.87
$ws.Range("F435").Value = 622
$ws.Range("G435").Value = 8365.9
$ws.Range("F436").Value = 519
$ws.Range("G436").Value = 6824.85
$ws.Range("F437").Value = 621
$ws.Range("G437").Value = 7955.01
$ws.Range("F438").Value = 284
$ws.Range("G438").Value = 7469.2
$ws.Range("F441").Value = 328
$ws.Range("G441").Value = 4201.68
$ws.Range("F443").Value = 445
$ws.Range("G443").Value = 2928.1
$ws.Range("F444").Value = 452
$ws.Range("G444").Value = 7331.44
$ws.Range("F445").Value = 130
$ws.Range("G445").Value = 2529.8
$ws.Range("F448").Value = 969
$ws.Range("G448").Value = 6288.81
$ws.Range("F450").Value = 355
$ws.Range("G450").Value = 9336.5
$ws.Range("B453").Value = 107777.11
$ws.Range("F560").Value = 31
$ws.Range("G560").Value = 1557.13
$ws.Range("B571").Value = 44500.94
$ws.Range("F577").Value = 96
$ws.Range("G577").Value = 2611.2
$ws.Range("F578").Value = 80
$ws.Range("G578").Value = 2176
$ws.Range("F579").Value = 30
$ws.Range("G579").Value = 816
$ws.Range("B580").Value = 65549.66
$ws.Range("F600").Value = 227
$ws.Range("G600").Value = 3600.22
$ws.Range("F605").Value = 86
$ws.Range("G605").Value = 2847.46
$ws.Range("B608").Value = 33083.53
$ws.Range("F637").Value = 0
$ws.Range("G637").Value = 0
$ws.Range("F640").Value = 2
$ws.Range("G640").Value = 266.64
$ws.Range("B642").Value = 3692.44
$ws.Range("F657").Value = 9
$ws.Range("G657").Value = 2121.66
$ws.Range("B664").Value = 22143.56
$ws.Range("F669").Value = 5
$ws.Range("G669").Value = 1079.45
$ws.Range("F670").Value = 211
$ws.Range("G670").Value = 14677.16
$ws.Range("F671").Value = 63
$ws.Range("G671").Value = 2611.35
$ws.Range("F672").Value = 13
$ws.Range("G672").Value = 698.49
$ws.Range("B676").Value = 45733.85
$ws.Range("F686").Value = 29
$ws.Range("G686").Value = 2279.11
$ws.Range("B694").Value = 24727.33
$ws.Range("F697").Value = 51
$ws.Range("G697").Value = 1907.4
$ws.Range("B702").Value = 10360.34
$ws.Range("F747").Value = 1126
$ws.Range("G747").Value = 183661.86
$ws.Range("F751").Value = 114
$ws.Range("G751").Value = 7695
$ws.Range("B752").Value = 215688.54
$ws.Range("B753").Value = 2294291.9
$ws.Range("B754").Value = 2294291.9
